# Updates the cryptos price/volume table (Sheet1) to the latest scraped
# snapshot. Most rows only get refreshed Price (D) and Volume/1h (E)
# figures; a couple of coins (InternetComputer/BitcoinCash and
# RenderToken/Stellar) swapped ranking order, so their whole rows
# (Coin, Link, Price, Volume) are rewritten.
#
# Price values that are plain decimal numbers (e.g. "115.30", "7.02")
# are prefixed with a leading apostrophe so Excel stores them as literal
# text instead of silently re-parsing them as numbers (which would trim
# trailing zeros / change precision). Values that already look like text
# to Excel (contain two dots, percent signs, letters, etc.) are set as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.788.54"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.290.95"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'115.30"
$ws.Range("E5").Value = "  +16.70%  "
$ws.Range("D6").Value = "'268.92"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("D10").Value = "'48.85"
$ws.Range("E10").Value = "  +8.27%  "
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").Value = "'8.96"
$ws.Range("E12").Value = "  +13.52%  "
$ws.Range("D13").Value = "'0.108"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "'15.82"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "2.634.30"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "'0.873"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").Value = "2.285.27"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "43.682.62"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'0.0000110"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "'7.02"
$ws.Range("E20").Value = "  +12.88%  "
$ws.Range("D21").Value = "'72.30"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'2.43"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'9.88"
$ws.Range("E23").Value = "  +8.54%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'233.20"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'2.94"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").Value = "'11.73"
$ws.Range("E26").Value = "  +4.21%  "
$ws.Range("D28").Value = "'42.91"
$ws.Range("E28").Value = "  +12.77%  "
$ws.Range("D29").Value = "'3.92"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "'173.67"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("D33").Value = "'0.0938"
$ws.Range("E33").Value = "  +5.30%  "
$ws.Range("D34").Value = "'21.61"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("E35").Value = "  +4.80%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'4.82"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.127"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("D39").Value = "'0.107"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").Value = "'3.84"
$ws.Range("E40").Value = "  +8.83%  "
$ws.Range("E41").Value = "  +20.08%  "
$ws.Range("E42").Value = "  +15.20%  "
$ws.Range("E43").Value = "  +3.57%  "
$ws.Range("E44").Value = "  +2.00%  "
$ws.Range("D45").Value = "'6.33"
$ws.Range("E45").Value = "  +21.51%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").Value = "'102.98"
$ws.Range("E49").Value = "  +4.69%  "
$ws.Range("E50").Value = "  +3.77%  "
$ws.Range("E51").Value = "  -2.05%  "
